# Apply the commit's row re-shuffle to the "Artfynd" sheet (rows 2-8).
#
# The edit re-orders the species observation rows and tweaks the
# "Taxonsorteringsordning" (column B) value for most of them, while row 1
# (headers) and the columns that don't vary between rows (C, I, P, S, T, U,
# V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-row data which must move together when rows are
# re-ordered.
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

# Snapshot the current ("before") values for rows 2-8 so the row move /
# permutation below can be applied without clobbering source data while we
# read it.
$before = @{}
for ($r = 2; $r -le 8; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $row
}

# New row order: key = destination row, value = source row (from $before)
# that its data (other than column B) originates from.
$rowSource = @{
    2 = 3
    3 = 4
    4 = 5
    5 = 8
    6 = 2
    7 = 6
    8 = 7
}

# New "Taxonsorteringsordning" (column B) values for each destination row.
$newB = @{
    2 = 56430
    3 = 78713
    4 = 77650
    5 = 77402
    6 = 96735
    7 = 78713
    8 = 95707
}

foreach ($destRow in 2..8) {
    $srcRow = $rowSource[$destRow]
    $srcData = $before[$srcRow]

    $ws.Range("A$destRow").Value = $srcData["A"]
    $ws.Range("B$destRow").Value = $newB[$destRow]
    $ws.Range("D$destRow").Value = $srcData["D"]
    $ws.Range("E$destRow").Value = $srcData["E"]
    $ws.Range("F$destRow").Value = $srcData["F"]
    $ws.Range("G$destRow").Value = $srcData["G"]
    $ws.Range("H$destRow").Value = $srcData["H"]
    $ws.Range("Q$destRow").Value = $srcData["Q"]
    $ws.Range("R$destRow").Value = $srcData["R"]

    $mVal = $srcData["M"]
    if ($mVal -eq $null -or $mVal -eq "") {
        $ws.Range("M$destRow").ClearContents()
    } else {
        $ws.Range("M$destRow").Value = $mVal
    }
}
